$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 31).Value = 2.840407848358154
$ws.Cells.Item(2, 32).Value = 12.974853515625
$ws.Cells.Item(3, 31).Value = 3.031013250350952
$ws.Cells.Item(3, 32).Value = 12.57294654846191
$ws.Cells.Item(4, 31).Value = 3.167018413543701
$ws.Cells.Item(4, 32).Value = 12.68288230895996
$ws.Cells.Item(5, 31).Value = 3.325519800186157
$ws.Cells.Item(5, 32).Value = 12.58668422698975
$ws.Cells.Item(6, 31).Value = 3.4575035572052
$ws.Cells.Item(6, 32).Value = 12.79450035095215
$ws.Cells.Item(7, 31).Value = 2.281842708587646
$ws.Cells.Item(7, 32).Value = 12.56695747375488
$ws.Cells.Item(8, 31).Value = 2.592242479324341
$ws.Cells.Item(8, 32).Value = 12.80287265777588
$ws.Cells.Item(9, 31).Value = 2.847051620483398
$ws.Cells.Item(9, 32).Value = 12.79043388366699
$ws.Cells.Item(10, 31).Value = 3.142376661300659
$ws.Cells.Item(10, 32).Value = 12.82507419586182
$ws.Cells.Item(11, 31).Value = 3.382583618164062
$ws.Cells.Item(11, 32).Value = 12.70583057403564
$ws.Cells.Item(12, 31).Value = 3.685203790664673
$ws.Cells.Item(12, 32).Value = 12.91453552246094
$ws.Cells.Item(13, 31).Value = 3.829819917678833
$ws.Cells.Item(13, 32).Value = 13.22307395935059
$ws.Cells.Item(14, 31).Value = 1.420221567153931
$ws.Cells.Item(14, 32).Value = 12.45138072967529
$ws.Cells.Item(15, 31).Value = 1.835132598876953
$ws.Cells.Item(15, 32).Value = 12.84426498413086
$ws.Cells.Item(16, 31).Value = 2.31755518913269
$ws.Cells.Item(16, 32).Value = 12.93618583679199
$ws.Cells.Item(17, 31).Value = 2.627319097518921
$ws.Cells.Item(17, 32).Value = 12.97330188751221
$ws.Cells.Item(18, 31).Value = 2.817667961120605
$ws.Cells.Item(18, 32).Value = 12.99754905700684
$ws.Cells.Item(19, 31).Value = 3.306522607803345
$ws.Cells.Item(19, 32).Value = 12.9640417098999
$ws.Cells.Item(20, 31).Value = 3.604533195495605
$ws.Cells.Item(20, 32).Value = 13.22229099273682
$ws.Cells.Item(21, 31).Value = 4.06130313873291
$ws.Cells.Item(21, 32).Value = 13.41179466247559
$ws.Cells.Item(22, 31).Value = 4.217324256896973
$ws.Cells.Item(22, 32).Value = 13.59045791625977
$ws.Cells.Item(23, 31).Value = 1.282296419143677
$ws.Cells.Item(23, 32).Value = 13.28883266448975
$ws.Cells.Item(24, 31).Value = 1.533604502677917
$ws.Cells.Item(24, 32).Value = 13.19435882568359
$ws.Cells.Item(25, 31).Value = 1.886880040168762
$ws.Cells.Item(25, 32).Value = 13.02475070953369
$ws.Cells.Item(26, 31).Value = 2.233947038650513
$ws.Cells.Item(26, 32).Value = 13.18631267547607
$ws.Cells.Item(27, 31).Value = 2.816362857818604
$ws.Cells.Item(27, 32).Value = 13.17830085754395
$ws.Cells.Item(28, 31).Value = 3.263370513916016
$ws.Cells.Item(28, 32).Value = 13.40875911712646
$ws.Cells.Item(29, 31).Value = 3.610644578933716
$ws.Cells.Item(29, 32).Value = 13.49584865570068
$ws.Cells.Item(30, 31).Value = 4.176106929779053
$ws.Cells.Item(30, 32).Value = 13.77007389068604
$ws.Cells.Item(31, 31).Value = 4.412979125976562
$ws.Cells.Item(31, 32).Value = 13.80388832092285
$ws.Cells.Item(32, 31).Value = 1.026296854019165
$ws.Cells.Item(32, 32).Value = 13.68306922912598
$ws.Cells.Item(33, 31).Value = 1.280116200447083
$ws.Cells.Item(33, 32).Value = 13.55428123474121
$ws.Cells.Item(34, 31).Value = 1.588118433952332
$ws.Cells.Item(34, 32).Value = 13.43269157409668
$ws.Cells.Item(35, 31).Value = 2.036976337432861
$ws.Cells.Item(35, 32).Value = 13.32980155944824
$ws.Cells.Item(36, 31).Value = 2.535278797149658
$ws.Cells.Item(36, 32).Value = 13.33436107635498
$ws.Cells.Item(37, 31).Value = 3.069841384887695
$ws.Cells.Item(37, 32).Value = 13.60497093200684
$ws.Cells.Item(38, 31).Value = 3.69943642616272
$ws.Cells.Item(38, 32).Value = 13.80571460723877
$ws.Cells.Item(39, 31).Value = 4.202789306640625
$ws.Cells.Item(39, 32).Value = 13.82440662384033
$ws.Cells.Item(40, 31).Value = 4.484659671783447
$ws.Cells.Item(40, 32).Value = 13.87829399108887
$ws.Cells.Item(41, 31).Value = 0.5137186646461487
$ws.Cells.Item(41, 32).Value = 13.91857719421387
$ws.Cells.Item(42, 31).Value = 1.213401317596436
$ws.Cells.Item(42, 32).Value = 13.89088344573975
$ws.Cells.Item(43, 31).Value = 1.401995062828064
$ws.Cells.Item(43, 32).Value = 13.74648761749268
$ws.Cells.Item(44, 31).Value = 2.001736879348755
$ws.Cells.Item(44, 32).Value = 13.81153774261475
$ws.Cells.Item(45, 31).Value = 2.496543169021606
$ws.Cells.Item(45, 32).Value = 13.73969841003418
$ws.Cells.Item(46, 31).Value = 3.128198862075806
$ws.Cells.Item(46, 32).Value = 13.81718349456787
$ws.Cells.Item(47, 31).Value = 3.537582159042358
$ws.Cells.Item(47, 32).Value = 13.97863960266113
$ws.Cells.Item(48, 31).Value = 4.129258632659912
$ws.Cells.Item(48, 32).Value = 14.03242778778076
$ws.Cells.Item(49, 31).Value = 4.395395278930664
$ws.Cells.Item(49, 32).Value = 13.9338321685791
$ws.Cells.Item(50, 31).Value = 1.137687921524048
$ws.Cells.Item(50, 32).Value = 14.14419078826904
$ws.Cells.Item(51, 31).Value = 1.28240966796875
$ws.Cells.Item(51, 32).Value = 13.869140625
$ws.Cells.Item(52, 31).Value = 1.502074241638184
$ws.Cells.Item(52, 32).Value = 14.07583332061768
$ws.Cells.Item(53, 31).Value = 2.02942967414856
$ws.Cells.Item(53, 32).Value = 14.04880714416504
$ws.Cells.Item(54, 31).Value = 2.383607387542725
$ws.Cells.Item(54, 32).Value = 13.95694065093994
$ws.Cells.Item(55, 31).Value = 2.923919916152954
$ws.Cells.Item(55, 32).Value = 13.96650886535645
$ws.Cells.Item(56, 31).Value = 3.538928985595703
$ws.Cells.Item(56, 32).Value = 14.06075477600098
$ws.Cells.Item(57, 31).Value = 4.00605297088623
$ws.Cells.Item(57, 32).Value = 14.11138916015625
$ws.Cells.Item(58, 31).Value = 4.343692779541016
$ws.Cells.Item(58, 32).Value = 13.91593360900879
$ws.Cells.Item(59, 31).Value = 1.568461537361145
$ws.Cells.Item(59, 32).Value = 14.43152809143066
$ws.Cells.Item(60, 31).Value = 1.712540864944458
$ws.Cells.Item(60, 32).Value = 14.28320407867432
$ws.Cells.Item(61, 31).Value = 1.909233689308167
$ws.Cells.Item(61, 32).Value = 14.4071159362793
$ws.Cells.Item(62, 31).Value = 2.268156290054321
$ws.Cells.Item(62, 32).Value = 14.25215816497803
$ws.Cells.Item(63, 31).Value = 2.72359561920166
$ws.Cells.Item(63, 32).Value = 14.18868160247803
$ws.Cells.Item(64, 31).Value = 3.324115991592407
$ws.Cells.Item(64, 32).Value = 14.13976860046387
$ws.Cells.Item(65, 31).Value = 3.923548221588135
$ws.Cells.Item(65, 32).Value = 14.05698585510254
$ws.Cells.Item(66, 31).Value = 1.746564507484436
$ws.Cells.Item(66, 32).Value = 14.54662990570068
$ws.Cells.Item(67, 31).Value = 2.034261226654053
$ws.Cells.Item(67, 32).Value = 14.33198738098145
$ws.Cells.Item(68, 31).Value = 2.308255910873413
$ws.Cells.Item(68, 32).Value = 14.19418621063232
$ws.Cells.Item(69, 31).Value = 2.591636180877686
$ws.Cells.Item(69, 32).Value = 14.37798118591309
$ws.Cells.Item(70, 31).Value = 3.16672158241272
$ws.Cells.Item(70, 32).Value = 14.07291889190674
